# Auto-generated Excel COM-interop script to apply the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows where only the Volume(1h) column (E) changes ---
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E33").Value = "  -6.65%  "
$ws.Range("E34").Value = "  -7.77%  "
$ws.Range("E42").Value = "  -0.04%  "

# --- Rows where Price (D) and Volume(1h) (E) change ---
$ws.Range("D2").Value = "61.930.01"
$ws.Range("E2").Value = "  -2.33%  "
$ws.Range("D3").Value = "2.497.99"
$ws.Range("E3").Value = "  -3.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.94"
$ws.Range("E5").Value = "  -3.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.90"
$ws.Range("E6").Value = "  -4.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.615"
$ws.Range("E8").Value = "  -1.16%  "
$ws.Range("D9").Value = "2.499.12"
$ws.Range("E9").Value = "  -3.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.107"
$ws.Range("E10").Value = "  -9.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.356"
$ws.Range("E13").Value = "  -6.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.14"
$ws.Range("E14").Value = "  -6.76%  "
$ws.Range("D15").Value = "2.950.67"
$ws.Range("E15").Value = "  -3.58%  "
$ws.Range("D16").Value = "61.862.78"
$ws.Range("E16").Value = "  -2.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000163"
$ws.Range("E17").Value = "  -8.15%  "
$ws.Range("D18").Value = "2.500.31"
$ws.Range("E18").Value = "  -3.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.09"
$ws.Range("E19").Value = "  -7.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.03"
$ws.Range("E20").Value = "  -5.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.19"
$ws.Range("E21").Value = "  -7.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "321.39"
$ws.Range("E22").Value = "  -5.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.91"
$ws.Range("E24").Value = "  -5.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.73"
$ws.Range("E25").Value = "  -3.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000101"
$ws.Range("E26").Value = "  -6.06%  "
$ws.Range("D27").Value = "2.622.59"
$ws.Range("E27").Value = "  -3.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.63"
$ws.Range("E32").Value = "  -3.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.57"
$ws.Range("E35").Value = "  -8.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.87"
$ws.Range("E36").Value = "  -9.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.87"
$ws.Range("E37").Value = "  -9.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.378"
$ws.Range("E39").Value = "  -5.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.50"
$ws.Range("E40").Value = "  -6.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "143.94"
$ws.Range("E41").Value = "  -6.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.69"
$ws.Range("E43").Value = "  -8.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.32"
$ws.Range("E44").Value = "  -2.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.29"
$ws.Range("E45").Value = "  -7.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "149.22"
$ws.Range("E46").Value = "  -4.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.58"
$ws.Range("E47").Value = "  -8.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.78"
$ws.Range("E48").Value = "  -10.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0535"
$ws.Range("E49").Value = "  -8.60%  "

# --- Rows with full Coin/Link/Price/Volume changes (re-ranked coins) ---
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.48"
$ws.Range("E29").Value = "  -4.73%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.38"
$ws.Range("E30").Value = "  -7.92%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "534.53"
$ws.Range("E31").Value = "  -7.06%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.589"
$ws.Range("E50").Value = "  -5.63%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0951"
$ws.Range("E51").Value = "  -4.90%  "

